$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are written with a leading apostrophe so Excel keeps them as
# literal text (e.g. "60.167.92", "0.0000109", "1.00", "  +5.78%  ")
# instead of coercing them into numbers and losing formatting/zeros.

# Row 2
$ws.Range("D2").Value = "'60.167.92"
$ws.Range("E2").Value = "'  +5.78%  "

# Row 3
$ws.Range("D3").Value = "'3.336.37"
$ws.Range("E3").Value = "'  +2.36%  "

# Row 4
$ws.Range("E4").Value = "'  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'411.99"
$ws.Range("E5").Value = "'  +3.83%  "

# Row 6
$ws.Range("D6").Value = "'111.23"
$ws.Range("E6").Value = "'  +1.60%  "

# Row 7
$ws.Range("D7").Value = "'0.584"
$ws.Range("E7").Value = "'  +4.26%  "

# Row 8
$ws.Range("E8").Value = "'  +0.16%  "

# Row 9
$ws.Range("D9").Value = "'0.633"
$ws.Range("E9").Value = "'  +1.07%  "

# Row 10
$ws.Range("D10").Value = "'39.44"
$ws.Range("E10").Value = "'  +0.40%  "

# Row 11
$ws.Range("D11").Value = "'0.0982"
$ws.Range("E11").Value = "'  +2.59%  "

# Row 12
$ws.Range("D12").Value = "'0.143"
$ws.Range("E12").Value = "'  +1.06%  "

# Row 13
$ws.Range("D13").Value = "'3.881.22"
$ws.Range("E13").Value = "'  +2.88%  "

# Row 14
$ws.Range("D14").Value = "'8.44"
$ws.Range("E14").Value = "'  +2.50%  "

# Row 15
$ws.Range("D15").Value = "'19.75"
$ws.Range("E15").Value = "'  +3.42%  "

# Row 16
$ws.Range("D16").Value = "'3.342.99"
$ws.Range("E16").Value = "'  +2.22%  "

# Row 17
$ws.Range("D17").Value = "'1.04"
$ws.Range("E17").Value = "'  +0.62%  "

# Row 18
$ws.Range("D18").Value = "'60.044.56"
$ws.Range("E18").Value = "'  +5.94%  "

# Row 19
$ws.Range("D19").Value = "'10.75"
$ws.Range("E19").Value = "'  -0.75%  "

# Row 20
$ws.Range("D20").Value = "'3.37"
$ws.Range("E20").Value = "'  +2.28%  "

# Row 21
$ws.Range("D21").Value = "'0.0000109"
$ws.Range("E21").Value = "'  +3.93%  "

# Row 22
$ws.Range("D22").Value = "'13.15"
$ws.Range("E22").Value = "'  +1.98%  "

# Row 23
$ws.Range("D23").Value = "'299.83"
$ws.Range("E23").Value = "'  -1.50%  "

# Row 24
$ws.Range("D24").Value = "'75.22"
$ws.Range("E24").Value = "'  +0.21%  "

# Row 25
$ws.Range("D25").Value = "'3.19"
$ws.Range("E25").Value = "'  +1.39%  "

# Row 26
$ws.Range("D26").Value = "'28.53"
$ws.Range("E26").Value = "'  +1.46%  "

# Row 27
$ws.Range("B27").Value = "'LEO"
$ws.Range("C27").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'4.48"
$ws.Range("E27").Value = "'  +2.16%  "

# Row 28
$ws.Range("B28").Value = "'Filecoin"
$ws.Range("C28").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'8.02"
$ws.Range("E28").Value = "'  +1.68%  "

# Row 29
$ws.Range("B29").Value = "'RenderToken"
$ws.Range("C29").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.70"
$ws.Range("E29").Value = "'  +6.40%  "

# Row 30
$ws.Range("B30").Value = "'Kaspa"
$ws.Range("C30").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.179"
$ws.Range("E30").Value = "'  +5.88%  "

# Row 31
$ws.Range("D31").Value = "'0.114"
$ws.Range("E31").Value = "'  +4.46%  "

# Row 32
$ws.Range("D32").Value = "'2.56"
$ws.Range("E32").Value = "'  +19.59%  "

# Row 33
$ws.Range("D33").Value = "'11.44"
$ws.Range("E33").Value = "'  +3.81%  "

# Row 34
$ws.Range("E34").Value = "'  +0.20%  "

# Row 35
$ws.Range("D35").Value = "'39.42"
$ws.Range("E35").Value = "'  +5.26%  "

# Row 36
$ws.Range("D36").Value = "'0.0504"
$ws.Range("E36").Value = "'  +5.03%  "

# Row 37
$ws.Range("D37").Value = "'52.28"
$ws.Range("E37").Value = "'  +1.53%  "

# Row 38
$ws.Range("B38").Value = "'FirstDigitalUSD"
$ws.Range("C38").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "'  +0.01%  "

# Row 39
$ws.Range("B39").Value = "'Stacks"
$ws.Range("C39").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.08"
$ws.Range("E39").Value = "'  -0.33%  "

# Row 40
$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "'  -3.17%  "

# Row 41
$ws.Range("D41").Value = "'137.67"
$ws.Range("E41").Value = "'  +1.79%  "

# Row 42
$ws.Range("B42").Value = "'Stellar"
$ws.Range("C42").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.123"
$ws.Range("E42").Value = "'  +2.25%  "

# Row 43
$ws.Range("B43").Value = "'TheGraph"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.292"
$ws.Range("E43").Value = "'  +3.90%  "

# Row 44
$ws.Range("D44").Value = "'1.91"
$ws.Range("E44").Value = "'  -0.42%  "

# Row 45
$ws.Range("D45").Value = "'3.92"
$ws.Range("E45").Value = "'  -2.25%  "

# Row 46
$ws.Range("D46").Value = "'16.84"
$ws.Range("E46").Value = "'  -1.80%  "

# Row 47
$ws.Range("D47").Value = "'2.27"
$ws.Range("E47").Value = "'  +8.97%  "

# Row 48
$ws.Range("D48").Value = "'22.40"
$ws.Range("E48").Value = "'  +1.60%  "

# Row 49
$ws.Range("D49").Value = "'2.194.69"
$ws.Range("E49").Value = "'  +2.25%  "

# Row 50
$ws.Range("E50").Value = "'  +1.28%  "

# Row 51
$ws.Range("D51").Value = "'1.99"
$ws.Range("E51").Value = "'  -0.81%  "
